$wb = $excel.ActiveWorkbook

# Sheets: 1 = Overview, 2 = zh-cn, 3 = de-de
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# The "0e26d396-746d-4b83-8f97-f8bcafd8b240" file's handback transform failed.
# Its Status (shared across Overview, zh-cn and de-de sheets) changes from
# "Ready for handoff" to "Handback transform failed".
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Populate the "Error Detail" column (L) for that row on the zh-cn and de-de
# sheets with the handback mismatch error message.
$wsZhCn.Range("L3").Value = "Handback file name: jssr2mv2.bj1 is different with handoff file name: 0e26d396-746d-4b83-8f97-f8bcafd8b240.89f604dbc6954c1823627dbb2a44a146f0750bb3.zh-cn."
$wsDeDe.Range("L3").Value = "Handback file name: jssr2mv2.bj1 is different with handoff file name: 0e26d396-746d-4b83-8f97-f8bcafd8b240.89f604dbc6954c1823627dbb2a44a146f0750bb3.de-de."
